$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2294.4707
$ws.Range("J19").Value = 2378.32
$ws.Range("L19").Value = 2378.32
$ws.Range("N19").Value = -2728.32
$ws.Range("H31").Value = 7834.75
$ws.Range("I31").Value = 5170
$ws.Range("K31").Value = 15510
$ws.Range("M31").Value = -15280
$ws.Range("H33").Value = 3844.45
$ws.Range("J33").Value = 3164
$ws.Range("L33").Value = 3164
$ws.Range("N33").Value = -3622
$ws.Range("H44").Value = 8000
$ws.Range("I44").Value = 8000
$ws.Range("K44").Value = 8000
$ws.Range("M44").Value = -7538
$ws.Range("H51").Value = 28038
$ws.Range("J51").Value = 32099.8
$ws.Range("L51").Value = 32099.8
$ws.Range("N51").Value = -33067.8
$ws.Range("H76").Value = 18246.25
$ws.Range("I76").Value = 17997.666
$ws.Range("J76").Value = 18992
$ws.Range("K76").Value = 17997.666
$ws.Range("L76").Value = 18992
$ws.Range("M76").Value = -17682.666
$ws.Range("N76").Value = -19622
$ws.Range("H79").Value = 18246.25
$ws.Range("I79").Value = 17997.666
$ws.Range("J79").Value = 18992
$ws.Range("K79").Value = 17997.666
$ws.Range("L79").Value = 18992
$ws.Range("M79").Value = -16905.666
$ws.Range("N79").Value = -21176
$ws.Range("H86").Value = 2823.5
$ws.Range("I86").Value = 2823.5
$ws.Range("K86").Value = 2823.5
$ws.Range("M86").Value = -1700.5
$ws.Range("H89").Value = 2823.5
$ws.Range("I89").Value = 2823.5
$ws.Range("K89").Value = 14117.5
$ws.Range("M89").Value = -8501.5
$ws.Range("H99").Value = 2881.2
$ws.Range("J99").Value = 6499.5
$ws.Range("L99").Value = 19498.5
$ws.Range("N99").Value = -22494.5
$ws.Range("H103").Value = 1178.1904
$ws.Range("I103").Value = 1016.6667
$ws.Range("J103").Value = 1205.1111
$ws.Range("K103").Value = 3050.0001
$ws.Range("L103").Value = 3615.3333
$ws.Range("M103").Value = -2464.0001
$ws.Range("N103").Value = -4787.3333
$ws.Range("H104").Value = 526.4
$ws.Range("I104").Value = 526.4
$ws.Range("K104").Value = 1579.2
$ws.Range("M104").Value = 167.8000000000002
$ws.Range("H111").Value = 9953.286
$ws.Range("I111").Value = 4338.8
$ws.Range("J111").Value = 23989.5
$ws.Range("K111").Value = 13016.4
$ws.Range("L111").Value = 71968.5
$ws.Range("M111").Value = -9949.400000000001
$ws.Range("N111").Value = -78102.5
$ws.Range("H113").Value = 2108.889
$ws.Range("I113").Value = 2546.25
$ws.Range("K113").Value = 2546.25
$ws.Range("M113").Value = 707.75
$ws.Range("H121").Value = 4441.5
$ws.Range("J121").Value = 4441.5
$ws.Range("L121").Value = 13324.5
$ws.Range("N121").Value = -16818.5
$ws.Range("H125").Value = 1448.2222
$ws.Range("I125").Value = 1247
$ws.Range("K125").Value = 11223
$ws.Range("M125").Value = -8763
$ws.Range("H129").Value = 4093.2856
$ws.Range("I129").Value = 2419.25
$ws.Range("K129").Value = 7257.75
$ws.Range("M129").Value = -2257.75
$ws.Range("H132").Value = 2563.75
$ws.Range("I132").Value = 2578.205
$ws.Range("K132").Value = 7734.615
$ws.Range("M132").Value = -5204.615
$ws.Range("H135").Value = 6557.0415
$ws.Range("I135").Value = 1576.2222
$ws.Range("K135").Value = 14185.9998
$ws.Range("M135").Value = -11650.9998
$ws.Range("H138").Value = 2559.4219
$ws.Range("I138").Value = 1358.75
$ws.Range("J138").Value = 3279.825
$ws.Range("K138").Value = 4076.25
$ws.Range("L138").Value = 9839.474999999999
$ws.Range("M138").Value = 1063.75
$ws.Range("N138").Value = -20119.475
$ws.Range("H141").Value = 2959.0667
$ws.Range("I141").Value = 2880.634
$ws.Range("K141").Value = 8641.902
$ws.Range("M141").Value = -3461.902

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").Value = $null
$ws.Range("H32").Value = 2152.3188
$ws.Range("I32").Value = 929.7705
$ws.Range("K32").Value = 929.7705
$ws.Range("M32").Value = -642.7705
$ws.Range("H45").Value = 2350.087
$ws.Range("I45").Value = 2107.2222
$ws.Range("K45").Value = 2107.2222
$ws.Range("M45").Value = -1730.2222
$ws.Range("H61").Value = 3932.3447
$ws.Range("I61").Value = 3160
$ws.Range("K61").Value = 3160
$ws.Range("M61").Value = -2948
$ws.Range("H63").Value = 2543.1428
$ws.Range("I63").Value = 2600.6667
$ws.Range("J63").Value = 2500
$ws.Range("K63").Value = 2600.6667
$ws.Range("L63").Value = 2500
$ws.Range("M63").Value = -1914.6667
$ws.Range("N63").Value = -3872
$ws.Range("H66").Value = 2543.1428
$ws.Range("I66").Value = 2600.6667
$ws.Range("J66").Value = 2500
$ws.Range("K66").Value = 13003.3335
$ws.Range("L66").Value = 12500
$ws.Range("M66").Value = -9571.333500000001
$ws.Range("N66").Value = -19364
$ws.Range("H74").Value = 2208.4
$ws.Range("I74").Value = 2166.7368
$ws.Range("K74").Value = 2166.7368
$ws.Range("M74").Value = -1292.7368
$ws.Range("H77").Value = 2208.4
$ws.Range("I77").Value = 2166.7368
$ws.Range("K77").Value = 10833.684
$ws.Range("M77").Value = -6465.684000000001
$ws.Range("H97").Value = 640.3
$ws.Range("I97").Value = 687.75
$ws.Range("J97").Value = 450.5
$ws.Range("K97").Value = 687.75
$ws.Range("L97").Value = 450.5
$ws.Range("M97").Value = -191.75
$ws.Range("N97").Value = -1442.5
$ws.Range("H122").Value = 3520.9524
$ws.Range("I122").Value = 3291.2354
$ws.Range("K122").Value = 9873.706200000001
$ws.Range("M122").Value = -7423.706200000001
$ws.Range("H132").Value = 2869.4138
$ws.Range("I132").Value = 2810.65
$ws.Range("K132").Value = 8431.950000000001
$ws.Range("M132").Value = -5901.950000000001
$ws.Range("H135").Value = 79881.19
$ws.Range("J135").Value = 79881.19
$ws.Range("L135").Value = 79881.19
$ws.Range("N135").Value = -90021.19
$ws.Range("H136").Value = 3932.3447
$ws.Range("I136").Value = 3160
$ws.Range("K136").Value = 9480
$ws.Range("M136").Value = -6930

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 14290376
$ws.Range("I86").Value = 33336398
$ws.Range("J86").Value = 5859.125
$ws.Range("K86").Value = 33336398
$ws.Range("L86").Value = 5859.125
$ws.Range("M86").Value = -33335275
$ws.Range("N86").Value = -8105.125
$ws.Range("H89").Value = 14290376
$ws.Range("I89").Value = 33336398
$ws.Range("J89").Value = 5859.125
$ws.Range("K89").Value = 166681990
$ws.Range("L89").Value = 29295.625
$ws.Range("M89").Value = -166676374
$ws.Range("N89").Value = -40527.625
$ws.Range("H94").Value = 2418.5
$ws.Range("I94").Value = 1388.7142
$ws.Range("J94").Value = 3619.9167
$ws.Range("K94").Value = 1388.7142
$ws.Range("L94").Value = 3619.9167
$ws.Range("M94").Value = -937.7141999999999
$ws.Range("N94").Value = -4521.9167
$ws.Range("H105").Value = 4877.391
$ws.Range("I105").Value = 4884
$ws.Range("J105").Value = 4833.3335
$ws.Range("K105").Value = 4884
$ws.Range("L105").Value = 4833.3335
$ws.Range("M105").Value = -3137
$ws.Range("N105").Value = -8327.333500000001
$ws.Range("H107").Value = 2156.7058
$ws.Range("I107").Value = 1872.6428
$ws.Range("K107").Value = 1872.6428
$ws.Range("M107").Value = 47.35719999999992

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 50000000
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").Value = $null
$ws.Range("H5").Value = 13602
$ws.Range("I5").Value = 10007
$ws.Range("J5").Value = 15399.5
$ws.Range("K5").Value = 10007
$ws.Range("L5").Value = 15399.5
$ws.Range("M5").Value = -9895
$ws.Range("N5").Value = -15623.5
$ws.Range("H16").Value = 11983.692
$ws.Range("I16").Value = 980.9091
$ws.Range("K16").Value = 980.9091
$ws.Range("M16").Value = -693.9091
$ws.Range("H31").Value = 3267.2373
$ws.Range("I31").Value = 5056.1055
$ws.Range("K31").Value = 5056.1055
$ws.Range("M31").Value = -4761.1055
$ws.Range("H34").Value = 3267.2373
$ws.Range("I34").Value = 5056.1055
$ws.Range("K34").Value = 5056.1055
$ws.Range("M34").Value = -4854.1055
$ws.Range("H62").Value = 2659.6
$ws.Range("I62").Value = 2749.5
$ws.Range("J62").Value = 2599.6667
$ws.Range("K62").Value = 2749.5
$ws.Range("L62").Value = 2599.6667
$ws.Range("M62").Value = -2125.5
$ws.Range("N62").Value = -3847.6667
$ws.Range("H65").Value = 2659.6
$ws.Range("I65").Value = 2749.5
$ws.Range("J65").Value = 2599.6667
$ws.Range("K65").Value = 13747.5
$ws.Range("L65").Value = 12998.3335
$ws.Range("M65").Value = -10627.5
$ws.Range("N65").Value = -19238.3335
$ws.Range("H70").Value = 49997.5
$ws.Range("J70").Value = 47663.332
$ws.Range("L70").Value = 47663.332
$ws.Range("N70").Value = -48293.332
$ws.Range("H73").Value = 49997.5
$ws.Range("J73").Value = 47663.332
$ws.Range("L73").Value = 47663.332
$ws.Range("N73").Value = -49847.332
$ws.Range("H86").Value = 18877.889
$ws.Range("I86").Value = 30600.25
$ws.Range("K86").Value = 30600.25
$ws.Range("M86").Value = -29477.25
$ws.Range("H89").Value = 18877.889
$ws.Range("I89").Value = 30600.25
$ws.Range("K89").Value = 153001.25
$ws.Range("M89").Value = -147385.25
$ws.Range("H99").Value = 7513.0713
$ws.Range("J99").Value = 10213.143
$ws.Range("L99").Value = 10213.143
$ws.Range("N99").Value = -13209.143
$ws.Range("H107").Value = 1015.6667
$ws.Range("J107").Value = 1023.75
$ws.Range("L107").Value = 1023.75
$ws.Range("N107").Value = -4863.75
$ws.Range("H113").Value = 11983.692
$ws.Range("I113").Value = 980.9091
$ws.Range("K113").Value = 980.9091
$ws.Range("M113").Value = 1189.0909
$ws.Range("H122").Value = 1198.75
$ws.Range("I122").Value = 1019.4
$ws.Range("J122").Value = 1497.6666
$ws.Range("K122").Value = 3058.2
$ws.Range("L122").Value = 4492.9998
$ws.Range("M122").Value = -608.1999999999998
$ws.Range("N122").Value = -9392.9998
$ws.Range("H126").Value = 7513.0713
$ws.Range("J126").Value = 10213.143
$ws.Range("L126").Value = 30639.429
$ws.Range("N126").Value = -35579.429
$ws.Range("H132").Value = 8337.556
$ws.Range("I132").Value = 8337.556
$ws.Range("K132").Value = 25012.668
$ws.Range("M132").Value = -22482.668

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4466.263
$ws.Range("I3").Value = 3580.0588
$ws.Range("K3").Value = 10740.1764
$ws.Range("M3").Value = -10628.1764
$ws.Range("H38").Value = 98.833336
$ws.Range("J38").Value = 93.25
$ws.Range("L38").Value = 279.75
$ws.Range("N38").Value = -973.75
$ws.Range("H60").Value = 716.2
$ws.Range("I60").Value = 395.25
$ws.Range("K60").Value = 1185.75
$ws.Range("M60").Value = -934.75
$ws.Range("H112").Value = 6900.515
$ws.Range("I112").Value = 5905.6665
$ws.Range("J112").Value = 7000
$ws.Range("K112").Value = 17716.9995
$ws.Range("L112").Value = 21000
$ws.Range("M112").Value = -16608.9995
$ws.Range("N112").Value = -23216
$ws.Range("H129").Value = 901.75
$ws.Range("J129").Value = 2000
$ws.Range("L129").Value = 6000
$ws.Range("N129").Value = -16000

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = $null
$ws.Range("H45").Value = 35588.445
$ws.Range("I45").Value = 80296
$ws.Range("K45").Value = 80296
$ws.Range("M45").Value = -79737
$ws.Range("H56").Value = 99
$ws.Range("I56").Value = 99
$ws.Range("K56").Value = 99
$ws.Range("M56").Value = 653
$ws.Range("H80").Value = 5227.357
$ws.Range("J80").Value = 5997.6665
$ws.Range("L80").Value = 5997.6665
$ws.Range("N80").Value = -7993.6665
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = $null
$ws.Range("H83").Value = 5227.357
$ws.Range("J83").Value = 5997.6665
$ws.Range("L83").Value = 29988.3325
$ws.Range("N83").Value = -39972.3325
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = $null
$ws.Range("H97").Value = 758.1
$ws.Range("I97").Value = 723
$ws.Range("K97").Value = 723
$ws.Range("M97").Value = -227
$ws.Range("H102").Value = 1620
$ws.Range("I102").Value = 1661.7142
$ws.Range("K102").Value = 1661.7142
$ws.Range("M102").Value = -39.71419999999989
$ws.Range("H113").Value = 19643.666
$ws.Range("I113").Value = 21099.125
$ws.Range("K113").Value = 21099.125
$ws.Range("M113").Value = -18929.125
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550
$ws.Range("H132").Value = 5987.528
$ws.Range("I132").Value = 6408.4136
$ws.Range("K132").Value = 19225.2408
$ws.Range("M132").Value = -16695.2408

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 6000
$ws.Range("J4").Value = 6000
$ws.Range("L4").Value = 6000
$ws.Range("N4").Value = -6226
$ws.Range("H7").Value = 18220.666
$ws.Range("I7").Value = 18220.666
$ws.Range("K7").Value = 18220.666
$ws.Range("M7").Value = -18108.666
$ws.Range("H16").Value = 1535.35
$ws.Range("I16").Value = 1587.3125
$ws.Range("J16").Value = 1327.5
$ws.Range("K16").Value = 1587.3125
$ws.Range("L16").Value = 1327.5
$ws.Range("M16").Value = -1417.3125
$ws.Range("N16").Value = -1667.5
$ws.Range("H28").Value = 6000
$ws.Range("J28").Value = 6000
$ws.Range("L28").Value = 6000
$ws.Range("N28").Value = -6464
$ws.Range("H37").Value = 6000
$ws.Range("J37").Value = 6000
$ws.Range("L37").Value = 6000
$ws.Range("N37").Value = -6214
$ws.Range("H40").Value = 2387.889
$ws.Range("I40").Value = 2387.889
$ws.Range("K40").Value = 2387.889
$ws.Range("M40").Value = -2251.889
$ws.Range("H46").Value = 3348.25
$ws.Range("I46").Value = 3531.3333
$ws.Range("K46").Value = 3531.3333
$ws.Range("M46").Value = -3343.3333
$ws.Range("H68").Value = 3315.6
$ws.Range("I68").Value = 2193
$ws.Range("K68").Value = 2193
$ws.Range("M68").Value = -1444
$ws.Range("H71").Value = 3315.6
$ws.Range("I71").Value = 2193
$ws.Range("K71").Value = 10965
$ws.Range("M71").Value = -7221
$ws.Range("H82").Value = 2895
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 2895
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 2895
$ws.Range("M82").Value = $null
$ws.Range("N82").Value = -3617
$ws.Range("H85").Value = 2895
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 2895
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 2895
$ws.Range("M85").Value = $null
$ws.Range("N85").Value = -5391
$ws.Range("H98").Value = 19999
$ws.Range("J98").Value = 19999
$ws.Range("L98").Value = 19999
$ws.Range("N98").Value = -25989
$ws.Range("H122").Value = 10276.529
$ws.Range("I122").Value = 12168.385
$ws.Range("J122").Value = 4128
$ws.Range("K122").Value = 36505.155
$ws.Range("L122").Value = 12384
$ws.Range("M122").Value = -34055.155
$ws.Range("N122").Value = -17284
$ws.Range("H126").Value = 18220.666
$ws.Range("I126").Value = 18220.666
$ws.Range("K126").Value = 54661.99800000001
$ws.Range("M126").Value = -52191.99800000001
$ws.Range("H141").Value = 88397.2
$ws.Range("J141").Value = 88774.664
$ws.Range("L141").Value = 88774.664
$ws.Range("N141").Value = -99134.664

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2059.1
$ws.Range("I113").Value = 1843.4445
$ws.Range("K113").Value = 5530.333500000001
$ws.Range("M113").Value = -3360.333500000001
$ws.Range("H122").Value = 1724.9
$ws.Range("I122").Value = 1626.12
$ws.Range("K122").Value = 4878.36
$ws.Range("M122").Value = -2428.36
$ws.Range("H127").Value = 99000
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").Value = $null
$ws.Range("H132").Value = 12967.4
$ws.Range("I132").Value = 10021.875
$ws.Range("K132").Value = 30065.625
$ws.Range("M132").Value = -27535.625
